$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 244
$ws.Range("I9").Value = 249.5
$ws.Range("J9").Value = 233
$ws.Range("K9").Value = 249.5
$ws.Range("L9").Value = 233
$ws.Range("M9").Value = -80.5
$ws.Range("N9").Value = -571

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 583.30304
$ws.Range("I28").Value = 210.7
$ws.Range("J28").Value = 1156.5385
$ws.Range("K28").Value = 210.7
$ws.Range("L28").Value = 1156.5385
$ws.Range("M28").Value = 274.3
$ws.Range("N28").Value = -2126.5385

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2540.2727
$ws.Range("I40").Value = 2642.2144
$ws.Range("K40").Value = 2642.2144
$ws.Range("M40").Value = -2467.2144

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1169.1
$ws.Range("I98").Value = 1198.3334
$ws.Range("J98").Value = 906
$ws.Range("K98").Value = 1198.3334
$ws.Range("L98").Value = 906
$ws.Range("M98").Value = 299.6666
$ws.Range("N98").Value = -3902

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 7105.3335
$ws.Range("I116").Value = 8126.6
$ws.Range("J116").Value = 1999
$ws.Range("K116").Value = 8126.6
$ws.Range("L116").Value = 1999
$ws.Range("M116").Value = -4684.6
$ws.Range("N116").Value = -8883

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1169.1
$ws.Range("I122").Value = 1198.3334
$ws.Range("J122").Value = 906
$ws.Range("K122").Value = 3595.0002
$ws.Range("L122").Value = 2718
$ws.Range("M122").Value = -1145.0002
$ws.Range("N122").Value = -7618

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 4314397.5
$ws.Range("I132").Value = 4906220
$ws.Range("K132").Value = 14718660
$ws.Range("M132").Value = -14716130

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1897.7858
$ws.Range("I137").Value = 1302.2858
$ws.Range("J137").Value = 3684.2856
$ws.Range("K137").Value = 3906.8574
$ws.Range("L137").Value = 11052.8568
$ws.Range("M137").Value = -1356.8574
$ws.Range("N137").Value = -16152.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 78967.92
$ws.Range("I2").Value = 2424.2856
$ws.Range("K2").Value = 2424.2856
$ws.Range("M2").Value = -2311.2856

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 24636.592
$ws.Range("I32").Value = 4712.923
$ws.Range("J32").Value = 142367.36
$ws.Range("K32").Value = 4712.923
$ws.Range("L32").Value = 142367.36
$ws.Range("M32").Value = -4425.923
$ws.Range("N32").Value = -142941.36

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 3725.75
$ws.Range("I88").Value = 4451.5
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 4451.5
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -4045.5
$ws.Range("N88").Value = -3812

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 3725.75
$ws.Range("I91").Value = 4451.5
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 4451.5
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -3047.5
$ws.Range("N91").Value = -5808

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 78967.92
$ws.Range("I116").Value = 2424.2856
$ws.Range("K116").Value = 2424.2856
$ws.Range("M116").Value = -130.2856000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 78967.92
$ws.Range("I3").Value = 2424.2856
$ws.Range("K3").Value = 2424.2856
$ws.Range("M3").Value = -2310.2856

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 38528.934
$ws.Range("I86").Value = 51513.684
$ws.Range("J86").Value = 2820.875
$ws.Range("K86").Value = 51513.684
$ws.Range("L86").Value = 2820.875
$ws.Range("M86").Value = -50390.684
$ws.Range("N86").Value = -5066.875

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 38528.934
$ws.Range("I89").Value = 51513.684
$ws.Range("J89").Value = 2820.875
$ws.Range("K89").Value = 257568.42
$ws.Range("L89").Value = 14104.375
$ws.Range("M89").Value = -251952.42
$ws.Range("N89").Value = -25336.375

# Row 111 (Leve Item ID 25789)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1995.254
$ws.Range("I134").Value = 2049.0408
$ws.Range("J134").Value = 1807
$ws.Range("K134").Value = 6147.1224
$ws.Range("L134").Value = 5421
$ws.Range("M134").Value = -3612.1224
$ws.Range("N134").Value = -10491

$ws = $wb.Worksheets.Item("CRP")
# Row 57 (Leve Item ID 3908)
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 11617.637
$ws.Range("J99").Value = 22162.8
$ws.Range("L99").Value = 22162.8
$ws.Range("N99").Value = -25158.8

# Row 112 (Leve Item ID 25796)
$ws.Range("H112").Value = 38000
$ws.Range("J112").Value = 38000
$ws.Range("L112").Value = 38000
$ws.Range("N112").Value = -40954

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 11617.637
$ws.Range("J126").Value = 22162.8
$ws.Range("L126").Value = 66488.39999999999
$ws.Range("N126").Value = -71428.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 14 (Leve Item ID 12886)
$ws.Range("H14").Value = 131.2
$ws.Range("I14").Value = 131.2
$ws.Range("K14").Value = 393.6
$ws.Range("M14").Value = -220.6

# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 673
$ws.Range("I23").Value = 284
$ws.Range("J23").Value = 889.1111
$ws.Range("K23").Value = 852
$ws.Range("L23").Value = 2667.3333
$ws.Range("M23").Value = -617
$ws.Range("N23").Value = -3137.3333

# Row 47 (Leve Item ID 4663)
$ws.Range("H47").Value = 135.46153
$ws.Range("I47").Value = 109.25
$ws.Range("J47").Value = 450
$ws.Range("K47").Value = 327.75
$ws.Range("L47").Value = 1350
$ws.Range("M47").Value = 103.25
$ws.Range("N47").Value = -2212

# Row 54 (Leve Item ID 4671)
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118

# Row 97 (Leve Item ID 19846)
$ws.Range("H97").Value = 700.6
$ws.Range("J97").Value = 750
$ws.Range("L97").Value = 2250
$ws.Range("N97").Value = -3242

# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 70717.94
$ws.Range("I98").Value = 1176.5
$ws.Range("J98").Value = 93898.414
$ws.Range("K98").Value = 3529.5
$ws.Range("L98").Value = 281695.242
$ws.Range("M98").Value = -2031.5
$ws.Range("N98").Value = -284691.242

# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 8746.223
$ws.Range("I121").Value = 9509.799999999999
$ws.Range("J121").Value = 8452.538
$ws.Range("K121").Value = 28529.4
$ws.Range("L121").Value = 25357.614
$ws.Range("M121").Value = -27219.4
$ws.Range("N121").Value = -27977.614

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 57773.816
$ws.Range("I70").Value = 104369.1
$ws.Range("J70").Value = 6001.278
$ws.Range("K70").Value = 104369.1
$ws.Range("L70").Value = 6001.278
$ws.Range("M70").Value = -104099.1
$ws.Range("N70").Value = -6541.278

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 57773.816
$ws.Range("I73").Value = 104369.1
$ws.Range("J73").Value = 6001.278
$ws.Range("K73").Value = 104369.1
$ws.Range("L73").Value = 6001.278
$ws.Range("M73").Value = -103433.1
$ws.Range("N73").Value = -7873.278

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 333335650
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 333335650
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984

$ws = $wb.Worksheets.Item("LTW")
# Row 111 (Leve Item ID 25820)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 868.875
$ws.Range("J113").Value = 1171.7142
$ws.Range("L113").Value = 3515.1426
$ws.Range("N113").Value = -7855.142599999999
